$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.1458753443493564
$ws.Range("H2").Value = 54.00547188224979
$ws.Range("I2").Value = -13.35872088791832
$ws.Range("G3").Value = 0.1234700994835908
$ws.Range("H3").Value = 86.81011056054332
$ws.Range("G4").Value = -0.02037819279444188
$ws.Range("H4").Value = -317.3611863864582
$ws.Range("G5").Value = -0.00068131118147376
$ws.Range("H5").Value = 94.91171971573858
$ws.Range("G6").Value = -0.2346069703862957
$ws.Range("H6").Value = -6.063084908930839
$ws.Range("G7").Value = -0.2152809577559405
$ws.Range("H7").Value = 13.84477856293245
$ws.Range("G8").Value = -0.3839107956671844
$ws.Range("H8").Value = -3.687970139870397
$ws.Range("G9").Value = -0.3408035595991667
$ws.Range("H9").Value = 14.50723130566895
$ws.Range("G10").Value = -0.03526918340668327
$ws.Range("H10").Value = -317.6453353147398
$ws.Range("G11").Value = 0.06993920524861373
$ws.Range("H11").Value = 535.3279148391908
$ws.Range("G12").Value = 0.2071601475528648
$ws.Range("H12").Value = -8.815863027497846
$ws.Range("G13").Value = 0.2443211867806952
$ws.Range("H13").Value = -7.222357388839186
$ws.Range("G14").Value = -0.04773665934998565
$ws.Range("H14").Value = -398.9578784797829
$ws.Range("G15").Value = -0.02417858951285413
$ws.Range("H15").Value = -219.7747488269223
$ws.Range("G16").Value = 0.1504293699359908
$ws.Range("H16").Value = 27.4834834277948
$ws.Range("G17").Value = 0.2109588251478912
$ws.Range("H17").Value = -3.598223334332349
$ws.Range("G18").Value = 0.03633657518354649
$ws.Range("H18").Value = -39.90637295203029
$ws.Range("G19").Value = 0.0932611936553083
$ws.Range("H19").Value = 3.520828114762166
$ws.Range("G20").Value = -0.1861103111122736
$ws.Range("H20").Value = -27.88364208494078
$ws.Range("G21").Value = -0.1941817681723289
$ws.Range("H21").Value = 2.826466804599221
$ws.Range("G22").Value = 0.03684709212262297
$ws.Range("H22").Value = -32.25173240534791
$ws.Range("G23").Value = 0.05051096004855238
$ws.Range("H23").Value = 23.67967795223643
$ws.Range("G24").Value = 0.1168416187481957
$ws.Range("H24").Value = 0.9567565352075369
$ws.Range("G25").Value = 0.1226324375084204
$ws.Range("H25").Value = -19.36109940629242
$ws.Range("G26").Value = -0.001287223466226757
$ws.Range("H26").Value = -102.4343829202037
$ws.Range("G27").Value = 0.05428923575006182
$ws.Range("H27").Value = 7.574143200579233
$ws.Range("G28").Value = 0.1642360865303584
$ws.Range("H28").Value = 7.40627469925303
$ws.Range("G29").Value = 0.1479946440853243
$ws.Range("H29").Value = -13.30288990318033
$ws.Range("G30").Value = -0.01357566026905447
$ws.Range("H30").Value = -169.3836246117841
$ws.Range("G31").Value = 0.04509623567294938
$ws.Range("H31").Value = 364.6715329841363
$ws.Range("G32").Value = 0.009226091716359876
$ws.Range("H32").Value = -75.26070750075917
$ws.Range("G33").Value = 0.02827864203931553
$ws.Range("H33").Value = 8.332183794959676
$ws.Range("G34").Value = 0.09599430385424816
$ws.Range("H34").Value = -24.98570159630874
$ws.Range("G35").Value = 0.1504382581715381
$ws.Range("H35").Value = 16.92644546615344
$ws.Range("G36").Value = -0.03377787299363369
$ws.Range("H36").Value = -324.7081284700537
$ws.Range("G37").Value = 0.01018000943508638
$ws.Range("H37").Value = -33.52693758138049
$ws.Range("G38").Value = -0.0407672325299857
$ws.Range("H38").Value = -1896.032182903549
$ws.Range("G39").Value = 0.01770532028529215
$ws.Range("H39").Value = 152.9953237884667
$ws.Range("G40").Value = 0.1274870066902942
$ws.Range("H40").Value = -13.59682983873568
$ws.Range("G41").Value = 0.1364195183795629
$ws.Range("H41").Value = -15.47695823927932
$ws.Range("G42").Value = 0.06980794863670302
$ws.Range("H42").Value = 8.120579600279795
$ws.Range("G43").Value = 0.0273683053257041
$ws.Range("H43").Value = -21.26612947237768
$ws.Range("G44").Value = 0.01697176416668398
$ws.Range("H44").Value = 20.25809908909864
$ws.Range("G45").Value = 0.03291857291657075
$ws.Range("H45").Value = -19.82429231400997
$ws.Range("G46").Value = -0.0812083982638054
$ws.Range("H46").Value = -23.3797721737125
$ws.Range("G47").Value = -0.03608251533009033
$ws.Range("H47").Value = 12.65365549679063
$ws.Range("G48").Value = -0.1367943811133881
$ws.Range("H48").Value = -8.588161562706555
$ws.Range("G49").Value = -0.1278903752515947
$ws.Range("H49").Value = 35.23928279500591
$ws.Range("G50").Value = 0.07875579745486189
$ws.Range("H50").Value = -27.6631222850219
$ws.Range("G51").Value = 0.1253540839503374
$ws.Range("H51").Value = 25.01550653904682
$ws.Range("G52").Value = 0.05870797280982552
$ws.Range("H52").Value = -1.528169954150575
$ws.Range("G53").Value = 0.04377433202443114
$ws.Range("H53").Value = -35.19011721707081
$ws.Range("G54").Value = -0.1174739124882164
$ws.Range("H54").Value = -68.0104940978466
$ws.Range("G55").Value = -0.05661924419879241
$ws.Range("H55").Value = 26.69338689281871
$ws.Range("G56").Value = 0.04920039590137342
$ws.Range("H56").Value = 7.363324067064066
$ws.Range("G57").Value = 0.1074112573784207
$ws.Range("H57").Value = 1977.475549460113
